$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Run 1)
$ws.Range("D8").Value = 224
$ws.Range("E8").Value = 113400
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = [double]"4.0443380850999999E-4"
$ws.Range("J8").Value = [double]"1.000016E-9"
$ws.Range("K8").Value = [double]"1.5028780220000001E-4"
$ws.Range("L8").Value = [double]"16.888791281"
$ws.Range("M8").Value = [double]"-2.8010999999999998E-12"
$ws.Range("N8").Value = [double]"-0.80774381747500001"
$ws.Range("O8").Value = [double]"-3.2082329999999999E-9"

# Row 9 (Run 2)
$ws.Range("D9").Value = 257
$ws.Range("E9").Value = 130032
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = [double]"4.0443380850199999E-4"
$ws.Range("J9").Value = [double]"1.0000000000000001E-9"
$ws.Range("K9").Value = [double]"1.5028780187200001E-4"
$ws.Range("L9").Value = [double]"16.88879127333"
$ws.Range("M9").Value = [double]"-4.8209999999999998E-12"
$ws.Range("N9").Value = [double]"-0.80774382132926403"
$ws.Range("O9").Value = [double]"-1.03109E-9"

# Row 10 (Run 3)
$ws.Range("D10").Value = 260
$ws.Range("E10").Value = 131544
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = [double]"4.0443380850700001E-4"
$ws.Range("J10").Value = [double]"1.000001E-9"
$ws.Range("K10").Value = [double]"1.5028780187399999E-4"
$ws.Range("L10").Value = [double]"16.888791273523101"
$ws.Range("M10").Value = [double]"-1.7320000000000001E-11"
$ws.Range("N10").Value = [double]"-0.80774382133240497"
$ws.Range("O10").Value = [double]"-1.0430119999999999E-9"

# Update selection to reflect the saved cursor position in the sheet view
$ws.Range("F12").Select()
